$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1975
$ws.Range("F6").Value = 1341
$ws.Range("F9").Value = 37
$ws.Range("F11").Value = 30
$ws.Range("F12").Value = 1859
$ws.Range("F14").Value = 1884
$ws.Range("F19").Value = 1615
$ws.Range("F21").Value = 27
$ws.Range("F23").Value = 2426
$ws.Range("F24").Value = 469
$ws.Range("F26").Value = 1058
$ws.Range("F27").Value = 4650
$ws.Range("F29").Value = 13
$ws.Range("F30").Value = 37
$ws.Range("F31").Value = 3
$ws.Range("F33").Value = 19
$ws.Range("F36").Value = 1251
$ws.Range("F37").Value = 5
$ws.Range("F38").Value = 1008
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 32
$ws.Range("F18").Value = 194
$ws.Range("F22").Value = 7
$ws.Range("F41").Value = 111
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2568
$ws.Range("F8").Value = 416
$ws.Range("F9").Value = 3134
$ws.Range("F10").Value = 637
$ws.Range("F11").Value = 909
$ws.Range("F12").Value = 328
$ws.Range("F13").Value = 48
$ws.Range("F14").Value = 74
$ws.Range("F15").Value = 19
$ws.Range("F16").Value = 328
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 416
$ws.Range("F7").Value = 3134
$ws.Range("F8").Value = 637
$ws.Range("F9").Value = 909
$ws.Range("F10").Value = 328
$ws.Range("F11").Value = 1975
$ws.Range("F13").Value = 48
$ws.Range("F14").Value = 48
$ws.Range("F15").Value = 74
$ws.Range("F16").Value = 74
$ws.Range("F17").Value = 30
$ws.Range("F18").Value = 19
$ws.Range("F19").Value = 1859
$ws.Range("F23").Value = 1884
$ws.Range("F26").Value = 1615
$ws.Range("F29").Value = 27
$ws.Range("F31").Value = 2426
$ws.Range("F32").Value = 469
$ws.Range("F34").Value = 1058
$ws.Range("F36").Value = 328
$ws.Range("F38").Value = 4650
$ws.Range("F41").Value = 37
$ws.Range("F45").Value = 111
$ws.Range("F48").Value = 1251
$ws.Range("F49").Value = 5
$ws.Range("F51").Value = 1008
